# 🔄 MAJ automatique BRVM via GitHub Actions
#
# Applies the daily BRVM data refresh to the "Recommandations" and
# "Top_YTD" sheets: updated variation figures, re-ranked rows (titles
# shuffled to their new sort position) and one newly appended row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Recommandations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Recommandations")

# Row 2 : BRVM - SERVICES PUBLICS
$ws.Range("D2").Value = 3306.62
$ws.Range("E2").Value = 109.45

# Row 3 : was CFAO MOTORS CI -> now SAFCA CI
$ws.Range("A3").Value = "SAFCA CI"
$ws.Range("D3").Value = 2680
$ws.Range("E3").Value = 675

# Row 4 : was SAFCA CI -> now CFAO MOTORS CI
$ws.Range("A4").Value = "CFAO MOTORS CI"
$ws.Range("D4").Value = 2665
$ws.Range("E4").Value = 645

# Row 5 : BRVM - AUTRES SECTEURS
$ws.Range("D5").Value = 2613.9
$ws.Range("E5").Value = 631.48

# Row 6 : NEI-CEDA CI
$ws.Range("D6").Value = 2365
$ws.Range("E6").Value = 595

# Row 7 : SETAO CI
$ws.Range("D7").Value = 2355

# Row 8 : UNIWAX CI
$ws.Range("D8").Value = 2310
$ws.Range("E8").Value = 585

# Row 9 : AIR LIQUIDE CI
$ws.Range("E9").Value = 520

# Row 10 : BRVM - DISTRIBUTION
$ws.Range("D10").Value = 1428.6
$ws.Range("E10").Value = 354.99

# Row 11 : BRVM - TRANSPORT
$ws.Range("D11").Value = 1386.64
$ws.Range("E11").Value = 347.58

# Row 12 : BRVM - AGRICULTURE
$ws.Range("D12").Value = 1233.07
$ws.Range("E12").Value = 309.14

# Row 13 : BRVM - INDUSTRIE
$ws.Range("D13").Value = 907.5599999999999
$ws.Range("E13").Value = 235.71

# Row 14 : BRVM - CONSOMMATION DE BASE
$ws.Range("D14").Value = 759.85
$ws.Range("E14").Value = 195.78

# Row 15 : BRVM-PRINCIPAL
$ws.Range("D15").Value = 726.61
$ws.Range("E15").Value = 183.74

# Row 16 : BRVM - INDUSTRIELS
$ws.Range("D16").Value = 545.72
$ws.Range("E16").Value = 139.71

# Row 17 : BRVM-PRESTIGE
$ws.Range("D17").Value = 519.2
$ws.Range("E17").Value = 130.07

# Row 18 : BRVM - FINANCES
$ws.Range("D18").Value = 491.56
$ws.Range("E18").Value = 123.48

# Row 19 : BRVM - SERVICES FINANCIERS
$ws.Range("D19").Value = 483.1
$ws.Range("E19").Value = 121.35

# Row 20 : BRVM - CONSOMMATION DISCRETIONNAIRE
$ws.Range("D20").Value = 425.1
$ws.Range("E20").Value = 104.92

# Row 21 : BRVM - ENERGIE
$ws.Range("D21").Value = 418.57
$ws.Range("E21").Value = 103.91

# Row 22 : BRVM - TELECOMMUNICATIONS
$ws.Range("D22").Value = 372.63
$ws.Range("E22").Value = 92.67

# Row 23 : UNILEVER CI (UNLC)
$ws.Range("D23").Value = 29.95

# Row 24 : FILTISAC CI (FTSC)
$ws.Range("B24").Value = 3
$ws.Range("D24").Value = 13.76
$ws.Range("E24").Value = 1.98
$ws.Range("F24").Value = "🟢 Achat"
$ws.Range("G24").Value = "✅ Renforcer"

# Row 25 : was CIE CI (CIEC) -> now BICI CI (BICC)
$ws.Range("A25").Value = "BICI CI (BICC)"
$ws.Range("D25").Value = 7.48
$ws.Range("E25").Value = 7.48

# Row 26 : was SODE CI (SDCC) -> now UNIWAX CI (UNXC)
$ws.Range("A26").Value = "UNIWAX CI (UNXC)"
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 6.56
$ws.Range("E26").Value = 7.27
$ws.Range("G26").Value = "👀 À surveiller"

# Row 29 : was SMB CI (SMBC) -> now CFAO MOTORS CI (CFAC)
$ws.Range("A29").Value = "CFAO MOTORS CI (CFAC)"
$ws.Range("D29").Value = 4.62
$ws.Range("E29").Value = 4.62

# Row 30 : was SOLIBRA CI (SLBC) -> now ECOBANK COTE D''IVOIRE (ECOC)
$ws.Range("A30").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws.Range("D30").Value = 4.46
$ws.Range("E30").Value = 4.46

# Row 31 : was BERNABE CI (BNBC) -> now SOLIBRA CI (SLBC)
$ws.Range("A31").Value = "SOLIBRA CI (SLBC)"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 3.83
$ws.Range("E31").Value = 3.83
$ws.Range("G31").Value = "➖ Neutre"

# Row 32 : was BANK OF AFRICA NG (BOAN) -> now SAFCA CI (SAFC)
$ws.Range("A32").Value = "SAFCA CI (SAFC)"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 2.99
$ws.Range("E32").Value = 2.99
$ws.Range("G32").Value = "➖ Neutre"

# Row 33 : was CFAO MOTORS CI (CFAC) -> now BANK OF AFRICA NG (BOAN)
$ws.Range("A33").Value = "BANK OF AFRICA NG (BOAN)"
$ws.Range("D33").Value = 1.94
$ws.Range("E33").Value = 5.1

# Row 36 : was UNIWAX CI (UNXC) -> now NESTLE CI (NTLC)
$ws.Range("A36").Value = "NESTLE CI (NTLC)"
$ws.Range("B36").Value = 0
$ws.Range("D36").Value = -0.85
$ws.Range("E36").Value = -0.85
$ws.Range("G36").Value = "➖ Neutre"

# Row 37 : was NESTLE CI (NTLC) -> now NEI-CEDA CI (NEIC)
$ws.Range("A37").Value = "NEI-CEDA CI (NEIC)"
$ws.Range("D37").Value = -1.68
$ws.Range("E37").Value = -1.68

# Row 38 : was ECOBANK COTE D''IVOIRE (ECOC) -> now BANK OF AFRICA ML (BOAM)
$ws.Range("A38").Value = "BANK OF AFRICA ML (BOAM)"
$ws.Range("B38").Value = 0
$ws.Range("D38").Value = -1.91
$ws.Range("E38").Value = -1.91
$ws.Range("G38").Value = "➖ Neutre"

# Row 39 : was NEI-CEDA CI (NEIC) -> now AIR LIQUIDE CI (SIVC)
$ws.Range("A39").Value = "AIR LIQUIDE CI (SIVC)"
$ws.Range("D39").Value = -1.92
$ws.Range("E39").Value = -1.92

# Row 40 : was BANK OF AFRICA ML (BOAM) -> now SERVAIR ABIDJAN CI (ABJC)
$ws.Range("A40").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws.Range("D40").Value = -2.42
$ws.Range("E40").Value = -2.42

# Row 41 : was VIVO ENERGY CI (SHEC) -> now BERNABE CI (BNBC)
$ws.Range("A41").Value = "BERNABE CI (BNBC)"
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = -2.45
$ws.Range("E41").Value = -2.24
$ws.Range("G41").Value = "👀 À surveiller"

# Row 42 : was SERVAIR ABIDJAN CI (ABJC) -> now NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws.Range("A42").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws.Range("D42").Value = -2.51
$ws.Range("E42").Value = -2.51

# Row 43 : was AIR LIQUIDE CI (SIVC) -> now BANK OF AFRICA BN (BOAB)
$ws.Range("A43").Value = "BANK OF AFRICA BN (BOAB)"
$ws.Range("D43").Value = -3.45
$ws.Range("E43").Value = -3.45

# Row 47 : was ONATEL BF (ONTBF) -> now LOTERIE NATIONALE DU BENIN (LNBB)
$ws.Range("A47").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = -4.26
$ws.Range("E47").Value = -4.26
$ws.Range("G47").Value = "➖ Neutre"

# Row 48 : TOTALENERGIES MARKETING SN (TTLS)
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = -4.47

# Row 49 (new) : ONATEL BF (ONTBF), appended at the end of the table
$ws.Range("A49").Value = "ONATEL BF (ONTBF)"
$ws.Range("B49").Value = 1
$ws.Range("C49").Value = 2
$ws.Range("D49").Value = -4.65
$ws.Range("E49").Value = -4.37
$ws.Range("F49").Value = "🟡 Observer"
$ws.Range("G49").Value = "👀 À surveiller"

# ---------------------------------------------------------------------
# Sheet 2: Top_YTD
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2 : BRVM - SERVICES PUBLICS
$ws2.Range("B2").Value = 8724110.869999999

# Row 3 : was CFAO MOTORS CI -> now SAFCA CI
$ws2.Range("A3").Value = "SAFCA CI"
$ws2.Range("B3").Value = 351265.62

# Row 4 : was SAFCA CI -> now CFAO MOTORS CI
$ws2.Range("A4").Value = "CFAO MOTORS CI"
$ws2.Range("B4").Value = 344419.66

# Row 5 : BRVM - AUTRES SECTEURS
$ws2.Range("B5").Value = 322021.59

# Row 6 : NEI-CEDA CI
$ws2.Range("B6").Value = 228213.76

# Row 7 : SETAO CI
$ws2.Range("B7").Value = 224928.67

# Row 8 : UNIWAX CI
$ws2.Range("B8").Value = 210488.18

# Row 10 : BRVM - DISTRIBUTION
$ws2.Range("B10").Value = 43569.23

# Row 11 : BRVM - TRANSPORT
$ws2.Range("B11").Value = 39701.3
